$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card6")

# New column M: "Event " header, matching the style used by the other
# header cells (bold font, thin border, centered horiz/top alignment).
$ws.Range("L1").Copy($ws.Range("M1"))
$ws.Range("M1").Value = "Event "

# The rest of the new column (rows 2-12) stays blank, just like the
# pandas/openpyxl export that produced the new column for every data row.
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 13).Style = "Normal"
}
